$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Shared string text run updates ---
$vol = $ws.Range("A8").Characters(21, 2)
$vol.Text = "37"

$d1 = $ws.Range("C9").Characters(27, 8)
$d1.Text = "9/11/2023"
$d2 = $ws.Range("C9").Characters(47, 9)
$d2.Text = "9/17/2023"

# --- Column width update (col H / 8) ---
$ws.Columns.Item(8).ColumnWidth = 8

# --- C27: convert placeholder text cell to a real number (style 15 / #,##0) ---
$ws.Range("C27").Value = 4
$ws.Range("C27").NumberFormat = "#,##0"

# --- Numeric cell value updates ---
$ws.Range("J15").Value = 18
$ws.Range("K15").Value = -44.444444444444
$ws.Range("C16").Value = 6
$ws.Range("D16").Value = 20
$ws.Range("E16").Value = -70
$ws.Range("F16").Value = 33
$ws.Range("G16").Value = 42
$ws.Range("H16").Value = -21.428571428571
$ws.Range("I16").Value = 339
$ws.Range("J16").Value = 440
$ws.Range("K16").Value = -22.954545454545
$ws.Range("L16").Value = 21.505376344086
$ws.Range("M16").Value = 208.181818181818
$ws.Range("N16").Value = -81.665765278528
$ws.Range("C17").Value = 10
$ws.Range("D17").Value = 12
$ws.Range("E17").Value = -16.666666666666
$ws.Range("F17").Value = 46
$ws.Range("G17").Value = 45
$ws.Range("H17").Value = 2.222222222222
$ws.Range("I17").Value = 362
$ws.Range("J17").Value = 322
$ws.Range("K17").Value = 12.422360248447
$ws.Range("L17").Value = 9.365558912386
$ws.Range("M17").Value = 158.571428571429
$ws.Range("N17").Value = -24.583333333333
$ws.Range("C18").Value = 8
$ws.Range("D18").Value = 10
$ws.Range("E18").Value = -20
$ws.Range("F18").Value = 27
$ws.Range("H18").Value = -48.076923076923
$ws.Range("I18").Value = 300
$ws.Range("J18").Value = 484
$ws.Range("K18").Value = -38.016528925619
$ws.Range("L18").Value = 2.040816326530
$ws.Range("M18").Value = 28.205128205128
$ws.Range("N18").Value = -84.407484407484
$ws.Range("C19").Value = 43
$ws.Range("D19").Value = 46
$ws.Range("E19").Value = -6.521739130434
$ws.Range("F19").Value = 159
$ws.Range("G19").Value = 190
$ws.Range("H19").Value = -16.315789473684
$ws.Range("I19").Value = 1654
$ws.Range("J19").Value = 1593
$ws.Range("K19").Value = 3.829252981795
$ws.Range("L19").Value = 81.161007667031
$ws.Range("M19").Value = 1.847290640394
$ws.Range("N19").Value = -75.687196824930
$ws.Range("F20").Value = 12
$ws.Range("H20").Value = 1100
$ws.Range("I20").Value = 52
$ws.Range("K20").Value = 8.333333333333
$ws.Range("L20").Value = 44.444444444444
$ws.Range("M20").Value = 205.882352941176
$ws.Range("N20").Value = -79.446640316205
$ws.Range("C21").Value = 70
$ws.Range("D21").Value = 89
$ws.Range("E21").Value = -21.348314606741
$ws.Range("F21").Value = 277
$ws.Range("G21").Value = 334
$ws.Range("H21").Value = -17.065868263473
$ws.Range("I21").Value = 2719
$ws.Range("J21").Value = 2911
$ws.Range("K21").Value = -6.595671590518
$ws.Range("L21").Value = 45.634708087841
$ws.Range("M21").Value = 27.413308341143
$ws.Range("N21").Value = -76.044052863436
$ws.Range("C22").Value = 4
$ws.Range("D22").Value = 7
$ws.Range("E22").Value = -42.857142857142
$ws.Range("F22").Value = 9
$ws.Range("G22").Value = 18
$ws.Range("H22").Value = -50
$ws.Range("I22").Value = 147
$ws.Range("J22").Value = 133
$ws.Range("K22").Value = 10.526315789473
$ws.Range("L22").Value = 38.679245283018
$ws.Range("M22").Value = 42.718446601941
$ws.Range("C24").Value = 100
$ws.Range("D24").Value = 74
$ws.Range("E24").Value = 35.135135135135
$ws.Range("F24").Value = 327
$ws.Range("G24").Value = 305
$ws.Range("H24").Value = 7.213114754098
$ws.Range("I24").Value = 2916
$ws.Range("J24").Value = 2344
$ws.Range("K24").Value = 24.402730375426
$ws.Range("L24").Value = 93.882978723404
$ws.Range("M24").Value = -12.955223880597
$ws.Range("C25").Value = 25
$ws.Range("D25").Value = 23
$ws.Range("E25").Value = 8.695652173913
$ws.Range("F25").Value = 95
$ws.Range("G25").Value = 69
$ws.Range("H25").Value = 37.681159420289
$ws.Range("I25").Value = 777
$ws.Range("J25").Value = 635
$ws.Range("K25").Value = 22.362204724409
$ws.Range("L25").Value = 23.529411764705
$ws.Range("M25").Value = 87.681159420289
$ws.Range("G26").Value = 4
$ws.Range("H26").Value = -75
$ws.Range("J26").Value = 23
$ws.Range("K26").Value = -17.391304347826
$ws.Range("D27").Value = 8
$ws.Range("E27").Value = -50
$ws.Range("F27").Value = 14
$ws.Range("H27").Value = -36.363636363636
$ws.Range("I27").Value = 160
$ws.Range("J27").Value = 166
$ws.Range("K27").Value = -3.614457831325
$ws.Range("L27").Value = 44.144144144144
$ws.Range("L30").Value = -67.857142857142
